$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "'23.09"
$ws.Range("D3").Style = "Normal"

$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D4").Value = "'5.399"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "3HuobiTokenHT"

$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D5").Value = "'0.05929"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "4CronosCRO"

$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'3.396"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "5GateTokenGT"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.8064"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "6MXTokenMX"

$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'0.9107"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "7FTXTokenFTT"

$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1415"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07431"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03331"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03055"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09326"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitMartTokenBMX"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'3.941"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001582"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04797"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005943"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16OneONE"

$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006099"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").Value = "'0.007493"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.004422"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009871"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20BitKanKAN"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.00007803"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.617"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("D24").Value = "'6.440"
$ws.Range("D24").Style = "Normal"

$ws.Range("D27").Value = "'0.1339"
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").Value = "'0.03876"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006206"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").Value = "'0.1069"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = "'0.002611"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.006438"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005186"
$ws.Range("D45").Style = "Normal"

$ws.Range("D48").Value = "'0.9805"
$ws.Range("D48").Style = "Normal"
